# Generate Report for Handoff
# Regenerates the localization-status report: the handoff file base name
# rolls from f004a19f-...-03759136a9af to 8ad3246c-...-126707fad4b0 (new
# markdown source + new xliff hash), and the handoff/generate timestamps
# advance a few seconds.

$wb = $excel.ActiveWorkbook

$oldGuid = "f004a19f-2db0-4156-9c37-03759136a9af"
$newGuid = "8ad3246c-fe12-4e35-bf27-126707fad4b0"

$newMdName  = $newGuid + ".md"
$newMdPath  = "e2e\" + $newMdName

$newZhXlf = $newGuid + ".803e9981aa1fac4744df834db6d1afd29301e6d9.zh-cn.xlf"
$newDeXlf = $newGuid + ".803e9981aa1fac4744df834db6d1afd29301e6d9.de-de.xlf"

$newGenerateDate = "2016-08-25 12:59:58"
$newZhHandoffDate = "2016-08-25 12:59:54"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = $newMdName
$wsOverview.Range("B2").Value2 = $newMdPath
$wsOverview.Range("G2").Value2 = $newGenerateDate
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdPath
}

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value2 = $newMdName
$wsZhCn.Range("G2").Value2 = $newZhXlf
$wsZhCn.Range("H2").Value2 = $newZhHandoffDate
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value2 = $newMdName
$wsDeDe.Range("G2").Value2 = $newDeXlf
$wsDeDe.Range("H2").Value2 = $newGenerateDate
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
